# Update the NATMI LR-pairs sheet with newly computed TPM-based statistics.
# Columns:
#   G/H  = Ligand average/total expression value   (depends on Sending cluster)
#   I/J  = Ligand derived specificity (avg/total)
#   M/N  = Receptor average/total expression value (depends on Target cluster)
#   O/P  = Receptor derived specificity (avg/total)
#   Q/R  = Edge average/total expression weight
#   S/T  = Edge average/total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ G=0.190922; H=0.572766; I=0.2296201981553925; J=0.2296201981553925;
            M=3.456265333333333; N=10.368796; O=0.009841535807677501; P=0.0098415358076775;
            Q=0.6598770899706666; R=5.938893809736; S=0.002259815402312299; T=0.002259815402312298 }
    3  = @{ G=0.190922; H=0.572766; I=0.2296201981553925; J=0.2296201981553925;
            O=0.8587907398420774; P=0.8587907398420773;
            Q=57.58210358373734; R=518.238932253636; S=0.197195699856554; T=0.1971956998565539 }
    4  = @{ G=0.190922; H=0.572766; I=0.2296201981553925; J=0.2296201981553925;
            O=0.1313677243502452; P=0.1313677243502452;
            Q=8.808234136859335; R=79.27410723173401; S=0.03016468289652629; T=0.03016468289652628 }
    5  = @{ I=0.5102601581298313; J=0.5102601581298313;
            M=3.456265333333333; N=10.368796; O=0.009841535807677501; P=0.0098415358076775;
            Q=1.466373563735111; R=13.197362073616; S=0.005021743617465919; T=0.005021743617465918 }
    6  = @{ I=0.5102601581298313; J=0.5102601581298313;
            O=0.8587907398420774; P=0.8587907398420773;
            S=0.4382066987122532; T=0.4382066987122532 }
    7  = @{ I=0.5102601581298313; J=0.5102601581298313;
            O=0.1313677243502452; P=0.1313677243502452;
            S=0.06703171580011219; T=0.06703171580011218 }
    8  = @{ I=0.2601196437147762; J=0.2601196437147761;
            M=3.456265333333333; N=10.368796; O=0.009841535807677501; P=0.0098415358076775;
            Q=0.7475256746471111; R=6.727731071824; S=0.002559976787899284; T=0.002559976787899283 }
    9  = @{ I=0.2601196437147762; J=0.2601196437147761;
            O=0.8587907398420774; P=0.8587907398420773;
            S=0.2233883412732702; T=0.2233883412732701 }
    10 = @{ I=0.2601196437147762; J=0.2601196437147761;
            O=0.1313677243502452; P=0.1313677243502452;
            S=0.03417132565360671; T=0.03417132565360669 }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
